$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.03073166666667
$ws.Range("H2").Value = 30.092195
$ws.Range("I2").Value = 0.5703560915538793
$ws.Range("J2").Value = 0.5703560915538795
$ws.Range("M2").Value = 108.845309
$ws.Range("N2").Value = 326.535927
$ws.Range("O2").Value = 0.3930063530400584
$ws.Range("P2").Value = 0.3930063530400583
$ws.Range("Q2").Value = 1091.798087754418
$ws.Range("R2").Value = 9826.182789789766
$ws.Range("S2").Value = 0.2241535674757718
$ws.Range("T2").Value = 0.2241535674757718

# Row 3
$ws.Range("G3").Value = 10.03073166666667
$ws.Range("H3").Value = 30.092195
$ws.Range("I3").Value = 0.5703560915538793
$ws.Range("J3").Value = 0.5703560915538795
$ws.Range("O3").Value = 0.4513455929560021
$ws.Range("P3").Value = 0.451345592956002
$ws.Range("Q3").Value = 1253.868420940053
$ws.Range("R3").Value = 11284.81578846048
$ws.Range("S3").Value = 0.2574277083384535
$ws.Range("T3").Value = 0.2574277083384535

# Row 4
$ws.Range("G4").Value = 10.03073166666667
$ws.Range("H4").Value = 30.092195
$ws.Range("I4").Value = 0.5703560915538793
$ws.Range("J4").Value = 0.5703560915538795
$ws.Range("M4").Value = 18.88973866666667
$ws.Range("N4").Value = 56.669216
$ws.Range("O4").Value = 0.06820493571538706
$ws.Range("P4").Value = 0.06820493571538705
$ws.Range("Q4").Value = 189.4778998187911
$ws.Range("R4").Value = 1705.30109836912
$ws.Range("S4").Value = 0.03890110055931176
$ws.Range("T4").Value = 0.03890110055931176

# Row 5
$ws.Range("G5").Value = 10.03073166666667
$ws.Range("H5").Value = 30.092195
$ws.Range("I5").Value = 0.5703560915538793
$ws.Range("J5").Value = 0.5703560915538795
$ws.Range("M5").Value = 13.37928666666667
$ws.Range("N5").Value = 40.13786
$ws.Range("O5").Value = 0.04830841776694433
$ws.Range("P5").Value = 0.04830841776694433
$ws.Range("Q5").Value = 134.2040344447445
$ws.Range("R5").Value = 1207.8363100027
$ws.Range("S5").Value = 0.02755300034670635
$ws.Range("T5").Value = 0.02755300034670636

# Row 6
$ws.Range("G6").Value = 10.03073166666667
$ws.Range("H6").Value = 30.092195
$ws.Range("I6").Value = 0.5703560915538793
$ws.Range("J6").Value = 0.5703560915538795
$ws.Range("M6").Value = 10.83857433333333
$ws.Range("N6").Value = 32.515723
$ws.Range("O6").Value = 0.03913470052160829
$ws.Range("P6").Value = 0.03913470052160829
$ws.Range("Q6").Value = 108.7188307868872
$ws.Range("R6").Value = 978.469477081985
$ws.Range("S6").Value = 0.02232071483363607
$ws.Range("T6").Value = 0.02232071483363607

# Row 7
$ws.Range("G7").Value = 7.556056333333333
$ws.Range("H7").Value = 22.668169
$ws.Range("I7").Value = 0.4296439084461207
$ws.Range("J7").Value = 0.4296439084461207
$ws.Range("M7").Value = 108.845309
$ws.Range("N7").Value = 326.535927
$ws.Range("O7").Value = 0.3930063530400584
$ws.Range("P7").Value = 0.3930063530400583
$ws.Range("Q7").Value = 822.4412864230736
$ws.Range("R7").Value = 7401.971577807663
$ws.Range("S7").Value = 0.1688527855642866
$ws.Range("T7").Value = 0.1688527855642866

# Row 8
$ws.Range("G8").Value = 7.556056333333333
$ws.Range("H8").Value = 22.668169
$ws.Range("I8").Value = 0.4296439084461207
$ws.Range("J8").Value = 0.4296439084461207
$ws.Range("O8").Value = 0.4513455929560021
$ws.Range("P8").Value = 0.451345592956002
$ws.Range("Q8").Value = 944.5273523460907
$ws.Range("R8").Value = 8500.746171114815
$ws.Range("S8").Value = 0.1939178846175486
$ws.Range("T8").Value = 0.1939178846175486

# Row 9
$ws.Range("G9").Value = 7.556056333333333
$ws.Range("H9").Value = 22.668169
$ws.Range("I9").Value = 0.4296439084461207
$ws.Range("J9").Value = 0.4296439084461207
$ws.Range("M9").Value = 18.88973866666667
$ws.Range("N9").Value = 56.669216
$ws.Range("O9").Value = 0.06820493571538706
$ws.Range("P9").Value = 0.06820493571538705
$ws.Range("Q9").Value = 142.7319294872782
$ws.Range("R9").Value = 1284.587365385504
$ws.Range("S9").Value = 0.02930383515607531
$ws.Range("T9").Value = 0.0293038351560753

# Row 10
$ws.Range("G10").Value = 7.556056333333333
$ws.Range("H10").Value = 22.668169
$ws.Range("I10").Value = 0.4296439084461207
$ws.Range("J10").Value = 0.4296439084461207
$ws.Range("M10").Value = 13.37928666666667
$ws.Range("N10").Value = 40.13786
$ws.Range("O10").Value = 0.04830841776694433
$ws.Range("P10").Value = 0.04830841776694433
$ws.Range("Q10").Value = 101.0946437531489
$ws.Range("R10").Value = 909.85179377834
$ws.Range("S10").Value = 0.02075541742023798
$ws.Range("T10").Value = 0.02075541742023798

# Row 11
$ws.Range("G11").Value = 7.556056333333333
$ws.Range("H11").Value = 22.668169
$ws.Range("I11").Value = 0.4296439084461207
$ws.Range("J11").Value = 0.4296439084461207
$ws.Range("M11").Value = 10.83857433333333
$ws.Range("N11").Value = 32.515723
$ws.Range("O11").Value = 0.03913470052160829
$ws.Range("P11").Value = 0.03913470052160829
$ws.Range("Q11").Value = 81.89687823568745
$ws.Range("R11").Value = 737.071904121187
$ws.Range("S11").Value = 0.01681398568797222
$ws.Range("T11").Value = 0.01681398568797222
